# GORES2021.xlsx update
# - Swap the two Coquimbo candidates in rows 20/21 on "tresquintos.cl"
# - Add four new Araucanía candidates (Luis Vivanco, Aucán Huilcamán,
#   César Vargas, Vicente Painel) to "tresquintos.cl"
# - Label the (until now unlabeled) status column on "todos" as "estado"
# - Update the remembered selections on both sheets

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "tresquintos.cl"
$ws2 = $wb.Worksheets.Item(2)   # "todos"

# --- Swap rows 20 and 21 (Marco Sulantay <-> Darío Molina) ---------------
$a20 = $ws1.Range("A20").Value
$e20 = $ws1.Range("E20").Value
$g20 = $ws1.Range("G20").Value

$a21 = $ws1.Range("A21").Value
$e21 = $ws1.Range("E21").Value
$g21 = $ws1.Range("G21").Value

$ws1.Range("A20").Value = $a21
$ws1.Range("E20").Value = $e21
$ws1.Range("G20").Value = $g21

$ws1.Range("A21").Value = $a20
$ws1.Range("E21").Value = $e20
$ws1.Range("G21").Value = $g20

# --- Insert the new Araucanía candidates ----------------------------------
# Before: row63=Luis Levi, row64=René Rubeska, row65=Eduardo Vicencio, row66=Eugenio Tuma

# Insert "Aucán Huilcamán" right after Luis Levi (row 64)
$ws1.Rows(64).Insert()
$ws1.Range("A64").Value = 83
$ws1.Range("B64").Value = 10
$ws1.Range("C64").Value = 9
$ws1.Range("D64").Value = "Araucanía"
$ws1.Range("E64").Value = "Aucán Huilcamán"
$ws1.Range("F64").Value = "Igualdad para Chile"
$ws1.Range("G64").Value = "IND"
$ws1.Range("H64").Value = 0

# Insert "Luis Vivanco" above the row just added, so it lands at row 64
$ws1.Rows(64).Insert()
$ws1.Range("A64").Value = 84
$ws1.Range("B64").Value = 10
$ws1.Range("C64").Value = 9
$ws1.Range("D64").Value = "Araucanía"
$ws1.Range("E64").Value = "Luis Vivanco"
$ws1.Range("F64").Value = "Humanicemos Chile"
$ws1.Range("G64").Value = "PH"
$ws1.Range("H64").Value = 0

# Insert "César Vargas" after Aucán Huilcamán, at row 66
$ws1.Rows(66).Insert()
$ws1.Range("A66").Value = 85
$ws1.Range("B66").Value = 10
$ws1.Range("C66").Value = 9
$ws1.Range("D66").Value = "Araucanía"
$ws1.Range("E66").Value = "César Vargas"
$ws1.Range("F66").Value = "Independientes Cristianos"
$ws1.Range("G66").Value = "IND"
$ws1.Range("H66").Value = 0

# Insert "Vicente Painel" after René Rubeska (now row 67), at row 68
$ws1.Rows(68).Insert()
$ws1.Range("A68").Value = 86
$ws1.Range("B68").Value = 10
$ws1.Range("C68").Value = 9
$ws1.Range("D68").Value = "Araucanía"
$ws1.Range("E68").Value = "Vicente Painel"
$ws1.Range("F68").Value = "Regionalistas Verdes"
$ws1.Range("G68").Value = "FRVS"
$ws1.Range("H68").Value = 0

# --- Label the status column on the "todos" sheet -------------------------
$ws2.Range("H1").Value = "estado"

# --- Restore remembered selections -----------------------------------------
$ws2.Activate()
$ws2.Range("H2").Select()

$ws1.Activate()
$ws1.Range("F11").Select()
